$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Revert "update task assign and deadlines":
# Unify deadline values in column E rows 22-48 to a single date string.
$ws.Range("E22:E48").Value = "8/4/2013 12:00AM"

# Restore the sheet view state (top-left cell + selection) as recorded before the
# "update task assign and deadlines" commit.
$ws.Range("F47").Select()
$ws.Application.ActiveWindow.ScrollRow = 32
